# Insert a new data row at row 237 (weekly price-update append), pushing the
# existing rows 237-286 down to 238-287. Excel's native row Insert takes care
# of shifting all the existing values/formatting down automatically, so we
# only need to populate the freshly inserted row with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(237).Insert()

$ws.Cells.Item(237, 1).Value  = 8
$ws.Cells.Item(237, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(237, 3).Value  = "Coquimbo"
$ws.Cells.Item(237, 4).Value  = 44641
$ws.Cells.Item(237, 5).Value  = 4
$ws.Cells.Item(237, 6).Value  = 100112032
$ws.Cells.Item(237, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(237, 8).Value  = "Sin especificar"
$ws.Cells.Item(237, 9).Value  = "Primera"
$ws.Cells.Item(237, 10).Value = 400
$ws.Cells.Item(237, 11).Value = 10000
$ws.Cells.Item(237, 12).Value = 11000
$ws.Cells.Item(237, 13).Value = 10500
$ws.Cells.Item(237, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(237, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(237, 16).Value = 175
$ws.Cells.Item(237, 17).Value = 60
$ws.Cells.Item(237, 18).Value = "Hortaliza"
